$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
$ws.Range("D2").Value = '58.762.14'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = '2.301.83'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''547.05'
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").Value = '''132.01'
$ws.Range("E6").Value = '  -2.70%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.572'
$ws.Range("E8").Value = '  -2.09%  '
$ws.Range("D9").Value = '2.300.31'
$ws.Range("E9").Value = '  -4.39%  '
$ws.Range("E10").Value = '  -2.60%  '
$ws.Range("D11").Value = '''5.51'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("D13").Value = '''0.333'
$ws.Range("E13").Value = '  -4.50%  '
$ws.Range("D14").Value = '''23.90'
$ws.Range("E14").Value = '  -2.94%  '
$ws.Range("D15").Value = '2.711.68'
$ws.Range("E15").Value = '  -4.42%  '
$ws.Range("D16").Value = '58.713.98'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("E17").Value = '  -2.94%  '
$ws.Range("D18").Value = '2.246.19'
$ws.Range("E18").Value = '  -6.54%  '
$ws.Range("E19").Value = '  -4.35%  '
$ws.Range("D20").Value = '''4.31'
$ws.Range("E20").Value = '  -4.42%  '
$ws.Range("D21").Value = '''314.53'
$ws.Range("E21").Value = '  -3.56%  '
$ws.Range("D22").Value = '''6.45'
$ws.Range("E22").Value = '  -4.12%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '''63.45'
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("D25").Value = '''0.168'
$ws.Range("E25").Value = '  -6.44%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -5.49%  '
$ws.Range("E28").Value = '  -5.89%  '
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("D30").Value = '''168.48'
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("D31").Value = '0.0₃0724'
$ws.Range("E31").Value = '  -5.13%  '
$ws.Range("E32").Value = '  +1.69%  '
$ws.Range("D33").Value = '''5.78'
$ws.Range("E33").Value = '  -5.44%  '
$ws.Range("E34").Value = '  -4.85%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '''17.79'
$ws.Range("E36").Value = '  -3.21%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '''1.25'
$ws.Range("E38").Value = '  -4.71%  '
$ws.Range("E39").Value = '  -5.12%  '
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("E41").Value = '  -4.90%  '
$ws.Range("D42").Value = '''298.18'
$ws.Range("E42").Value = '  -7.22%  '
$ws.Range("D43").Value = '''140.65'
$ws.Range("E43").Value = '  -3.37%  '
$ws.Range("D44").Value = '''3.43'
$ws.Range("E44").Value = '  -4.38%  '
$ws.Range("D45").Value = '''0.0952'
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("D46").Value = '''0.0500'
$ws.Range("E46").Value = '  -2.43%  '
$ws.Range("D47").Value = '''0.555'
$ws.Range("E47").Value = '  -3.36%  '
$ws.Range("D48").Value = '''18.42'
$ws.Range("E48").Value = '  -7.13%  '
$ws.Range("E49").Value = '  -2.78%  '
$ws.Range("E50").Value = '  -3.50%  '
$ws.Range("E51").Value = '  -0.29%  '
